# The workbook's first sheet ("Program Overview") had a blank separator
# row (row 9) between the FY2012 row and the FY2013-2016 rows. This edit
# removes that blank row, which shifts the FY2013-2016 data rows up by one
# (old rows 10-13 become new rows 9-12), shrinking the used range from
# A1:F13 to A1:F12 and moving the active selection to D9.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Program Overview")

# Delete the blank row 9 - this shifts everything below it up by one row.
$ws.Rows.Item(9).Delete()

# Update the active selection to match the post-edit state (D9).
$ws.Range("D9").Select()
